# Atualização de bases das ligas, do dia: 11-04-2024 às 00:31
#
# This edit swaps the data (all columns except A/Id, C/Div, D/Div Original
# Name and E/Date, which stay fixed per physical row) between several pairs
# of adjacent match rows. In other words, for each pair the two matches
# traded places while the row-local "id" sequence (column A) and the
# league/date metadata (columns C, D, E) were left untouched.
#
# Row pairs (1-based worksheet rows) that are swapped:
#   18 <-> 19
#   59 <-> 60
#   63 <-> 64
#   130 <-> 132
#   134 <-> 135
#   145 <-> 146

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together with each match record.
# (A = id, C = Div, D = Div Original Name, E = Date are intentionally excluded
#  because they remain identical / unchanged for the swapped rows.)
$swapCols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

$rowPairs = @(
    , @(18, 19)
    , @(59, 60)
    , @(63, 64)
    , @(130, 132)
    , @(134, 135)
    , @(145, 146)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $swapCols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $cell1 = $ws.Range($addr1)
        $cell2 = $ws.Range($addr2)

        $val1 = $cell1.Value2
        $val2 = $cell2.Value2

        $cell1.Value2 = $val2
        $cell2.Value2 = $val1
    }
}
